$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2250
$ws.Range("J29").Value = 2600
$ws.Range("L29").Value = 7800
$ws.Range("N29").Value = -8362

$ws.Range("H33").Value = 99.64706
$ws.Range("I33").Value = 79.888885
$ws.Range("J33").Value = 121.875
$ws.Range("K33").Value = 79.888885
$ws.Range("L33").Value = 121.875
$ws.Range("M33").Value = 149.111115
$ws.Range("N33").Value = -579.875

$ws.Range("H37").Value = 2208
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 2208
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 6624
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -6876

$ws.Range("H107").Value = 880.2727
$ws.Range("I107").Value = 522.8125
$ws.Range("J107").Value = 1833.5
$ws.Range("K107").Value = 522.8125
$ws.Range("L107").Value = 1833.5
$ws.Range("M107").Value = 1397.1875
$ws.Range("N107").Value = -5673.5

$ws.Range("H112").Value = 7455.3335
$ws.Range("J112").Value = 7455.3335
$ws.Range("L112").Value = 22366.0005
$ws.Range("N112").Value = -24582.0005

$ws.Range("H121").Value = 1492.5454
$ws.Range("J121").Value = 1696.4445
$ws.Range("L121").Value = 5089.333500000001
$ws.Range("N121").Value = -8583.333500000001

$ws.Range("H137").Value = 2284.7144
$ws.Range("I137").Value = 1998.6
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 5995.799999999999
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -3445.799999999999
$ws.Range("N137").Value = -14100

$ws.Range("H138").Value = 6569.6377
$ws.Range("I138").Value = 5415
$ws.Range("J138").Value = 6870.8477
$ws.Range("K138").Value = 16245
$ws.Range("L138").Value = 20612.5431
$ws.Range("M138").Value = -11105
$ws.Range("N138").Value = -30892.5431

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3625.3135
$ws.Range("I32").Value = 2834.2932
$ws.Range("K32").Value = 2834.2932
$ws.Range("M32").Value = -2547.2932

$ws.Range("H61").Value = 3088.2354
$ws.Range("I61").Value = 1958.4166
$ws.Range("J61").Value = 5799.8
$ws.Range("K61").Value = 1958.4166
$ws.Range("L61").Value = 5799.8
$ws.Range("M61").Value = -1746.4166
$ws.Range("N61").Value = -6223.8

$ws.Range("H94").Value = 99996
$ws.Range("J94").Value = 99996
$ws.Range("L94").Value = 99996
$ws.Range("N94").Value = -101798

$ws.Range("H110").Value = 1111
$ws.Range("I110").Value = 1111
$ws.Range("K110").Value = 1111
$ws.Range("M110").Value = 934

$ws.Range("H132").Value = 2515.913
$ws.Range("I132").Value = 2542.5454
$ws.Range("J132").Value = 2491.5
$ws.Range("K132").Value = 7627.6362
$ws.Range("L132").Value = 7474.5
$ws.Range("M132").Value = -5097.6362
$ws.Range("N132").Value = -12534.5

$ws.Range("H136").Value = 3088.2354
$ws.Range("I136").Value = 1958.4166
$ws.Range("J136").Value = 5799.8
$ws.Range("K136").Value = 5875.2498
$ws.Range("L136").Value = 17399.4
$ws.Range("M136").Value = -3325.2498
$ws.Range("N136").Value = -22499.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 6066.1763
$ws.Range("I80").Value = 26
$ws.Range("J80").Value = 7360.5
$ws.Range("K80").Value = 26
$ws.Range("L80").Value = 7360.5
$ws.Range("M80").Value = 972
$ws.Range("N80").Value = -9356.5

$ws.Range("H83").Value = 6066.1763
$ws.Range("I83").Value = 26
$ws.Range("J83").Value = 7360.5
$ws.Range("K83").Value = 130
$ws.Range("L83").Value = 36802.5
$ws.Range("M83").Value = 4862
$ws.Range("N83").Value = -46786.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1413.569
$ws.Range("I31").Value = 1016.7083
$ws.Range("J31").Value = 1693.7059
$ws.Range("K31").Value = 1016.7083
$ws.Range("L31").Value = 1693.7059
$ws.Range("M31").Value = -721.7083
$ws.Range("N31").Value = -2283.7059

$ws.Range("H34").Value = 1413.569
$ws.Range("I34").Value = 1016.7083
$ws.Range("J34").Value = 1693.7059
$ws.Range("K34").Value = 1016.7083
$ws.Range("L34").Value = 1693.7059
$ws.Range("M34").Value = -814.7083
$ws.Range("N34").Value = -2097.7059

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3356.9622
$ws.Range("I68").Value = 1402.4
$ws.Range("J68").Value = 4128.5
$ws.Range("K68").Value = 4207.200000000001
$ws.Range("L68").Value = 12385.5
$ws.Range("M68").Value = -3396.200000000001
$ws.Range("N68").Value = -14007.5

$ws.Range("H71").Value = 3356.9622
$ws.Range("I71").Value = 1402.4
$ws.Range("J71").Value = 4128.5
$ws.Range("K71").Value = 12621.6
$ws.Range("L71").Value = 37156.5
$ws.Range("M71").Value = -8565.6
$ws.Range("N71").Value = -45268.5

$ws.Range("H112").Value = 3499.5
$ws.Range("I112").Value = 1999
$ws.Range("J112").Value = 5000
$ws.Range("K112").Value = 5997
$ws.Range("L112").Value = 15000
$ws.Range("M112").Value = -4889
$ws.Range("N112").Value = -17216

$ws.Range("H114").Value = 9524932
$ws.Range("I114").Value = 545
$ws.Range("J114").Value = 17858770
$ws.Range("K114").Value = 1635
$ws.Range("L114").Value = 53576310
$ws.Range("M114").Value = 1619
$ws.Range("N114").Value = -53582818

$ws.Range("H129").Value = 32374
$ws.Range("I129").Value = 536.2222
$ws.Range("J129").Value = 52841.145
$ws.Range("K129").Value = 1608.6666
$ws.Range("L129").Value = 158523.435
$ws.Range("M129").Value = 3391.3334
$ws.Range("N129").Value = -168523.435

$ws.Range("H131").Value = 21770958
$ws.Range("J131").Value = 55912.46
$ws.Range("L131").Value = 167737.38
$ws.Range("N131").Value = -177817.38

$ws.Range("H137").Value = 4306.778
$ws.Range("I137").Value = 1288.1666
$ws.Range("J137").Value = 5816.0835
$ws.Range("K137").Value = 3864.4998
$ws.Range("L137").Value = 17448.2505
$ws.Range("M137").Value = 1235.5002
$ws.Range("N137").Value = -27648.2505

$ws.Range("H140").Value = 3012
$ws.Range("I140").Value = 1669.8
$ws.Range("J140").Value = 4130.5
$ws.Range("K140").Value = 5009.4
$ws.Range("L140").Value = 12391.5
$ws.Range("M140").Value = 170.6000000000004
$ws.Range("N140").Value = -22751.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4290.95
$ws.Range("I70").Value = 4509.923
$ws.Range("J70").Value = 3884.2856
$ws.Range("K70").Value = 4509.923
$ws.Range("L70").Value = 3884.2856
$ws.Range("M70").Value = -4239.923
$ws.Range("N70").Value = -4424.2856

$ws.Range("H73").Value = 4290.95
$ws.Range("I73").Value = 4509.923
$ws.Range("J73").Value = 3884.2856
$ws.Range("K73").Value = 4509.923
$ws.Range("L73").Value = 3884.2856
$ws.Range("M73").Value = -3573.923
$ws.Range("N73").Value = -5756.2856

$ws.Range("H93").Value = 27125
$ws.Range("J93").Value = 27125
$ws.Range("L93").Value = 27125
$ws.Range("N93").Value = -30869

$ws.Range("H132").Value = 2407840.8
$ws.Range("I132").Value = 6412816
$ws.Range("J132").Value = 4855.6
$ws.Range("K132").Value = 19238448
$ws.Range("L132").Value = 14566.8
$ws.Range("M132").Value = -19235918
$ws.Range("N132").Value = -19626.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 5374.4443
$ws.Range("I32").Value = 4103.2856
$ws.Range("J32").Value = 9823.5
$ws.Range("K32").Value = 4103.2856
$ws.Range("L32").Value = 9823.5
$ws.Range("M32").Value = -3786.2856
$ws.Range("N32").Value = -10457.5

$ws.Range("H132").Value = 2825.8367
$ws.Range("I132").Value = 1404.6957
$ws.Range("J132").Value = 4083
$ws.Range("K132").Value = 4214.0871
$ws.Range("L132").Value = 12249
$ws.Range("M132").Value = -1684.0871
$ws.Range("N132").Value = -17309

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1945
$ws.Range("I62").Value = 1945
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1945
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1321
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 1945
$ws.Range("I65").Value = 1945
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 9725
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -6605
$ws.Range("N65").ClearContents()

$ws.Range("H132").Value = 1371.3636
$ws.Range("I132").Value = 984.86664
$ws.Range("J132").Value = 2199.5715
$ws.Range("K132").Value = 2954.59992
$ws.Range("L132").Value = 6598.7145
$ws.Range("M132").Value = -424.5999199999997
$ws.Range("N132").Value = -11658.7145
